$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2014/12)
$ws.Range("D2").Value = 2167
$ws.Range("E2").Value = -152
$ws.Range("F2").Value = -152
$ws.Range("G2").Value = 16
$ws.Range("H2").Value = -2
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 1942
$ws.Range("L2").Value = 552
$ws.Range("M2").Value = 1390
$ws.Range("N2").Value = 1387
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 69
$ws.Range("Q2").Value = -111
$ws.Range("R2").Value = 337
$ws.Range("S2").Value = -236
$ws.Range("T2").Value = 12
$ws.Range("U2").Value = -123
$ws.Range("V2").Value = 280
$ws.Range("W2").Value = -7.03
$ws.Range("X2").Value = -0.08
$ws.Range("Y2").Value = 0.12
$ws.Range("Z2").Value = -0.08
$ws.Range("AA2").Value = 39.69
$ws.Range("AB2").Value = 1969
$ws.Range("AC2").Value = 25
$ws.Range("AD2").Value = 334.14
$ws.Range("AE2").Value = 21410
$ws.Range("AF2").Value = 0.39
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 2.41
$ws.Range("AI2").Value = 500.3
$ws.Range("AJ2").Value = 6867945

# Row 3 (2015/12)
$ws.Range("D3").Value = 2092
$ws.Range("E3").Value = -82
$ws.Range("F3").Value = -76
$ws.Range("G3").Value = 98
$ws.Range("H3").Value = 62
$ws.Range("I3").Value = 65
$ws.Range("J3").Value = -3
$ws.Range("K3").Value = 1694
$ws.Range("L3").Value = 346
$ws.Range("M3").Value = 1348
$ws.Range("N3").Value = 1348
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 69
$ws.Range("Q3").Value = -75
$ws.Range("R3").Value = 269
$ws.Range("S3").Value = -201
$ws.Range("T3").Value = 11
$ws.Range("U3").Value = -86
$ws.Range("V3").Value = 93
$ws.Range("W3").Value = -3.9
$ws.Range("X3").Value = 2.97
$ws.Range("Y3").Value = 4.76
$ws.Range("Z3").Value = 3.42
$ws.Range("AA3").Value = 25.65
$ws.Range("AB3").Value = 1920.51
$ws.Range("AC3").Value = 949
$ws.Range("AD3").Value = 10.65
$ws.Range("AE3").Value = 20812
$ws.Range("AF3").Value = 0.49
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 1.98
$ws.Range("AI3").Value = 19.88
$ws.Range("AJ3").Value = 6867945

# Row 4 (2016/12)
$ws.Range("D4").Value = 2075
$ws.Range("E4").Value = -20
$ws.Range("F4").Value = -20
$ws.Range("G4").Value = -24
$ws.Range("H4").Value = -23
$ws.Range("I4").Value = -22
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 1644
$ws.Range("L4").Value = 343
$ws.Range("M4").Value = 1302
$ws.Range("N4").Value = 1316
$ws.Range("O4").Value = -14
$ws.Range("P4").Value = 69
$ws.Range("Q4").Value = -29
$ws.Range("R4").Value = 50
$ws.Range("S4").Value = -2
$ws.Range("T4").Value = 5
$ws.Range("U4").Value = -35
$ws.Range("V4").Value = 95
$ws.Range("W4").Value = -0.95
$ws.Range("X4").Value = -1.1
$ws.Range("Y4").Value = -1.67
$ws.Range("Z4").Value = -1.37
$ws.Range("AA4").Value = 26.31
$ws.Range("AB4").Value = 1861
$ws.Range("AC4").Value = -323
$ws.Range("AD4").Value = -27.38
$ws.Range("AE4").Value = 20317
$ws.Range("AF4").Value = 0.44
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 2.26
$ws.Range("AI4").Value = -58.36
$ws.Range("AJ4").Value = 6867945

# Row 5 (2017/12)
$ws.Range("D5").Value = 2094
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = -11
$ws.Range("H5").Value = 16
$ws.Range("I5").Value = 11
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 1637
$ws.Range("L5").Value = 343
$ws.Range("M5").Value = 1294
$ws.Range("N5").Value = 1304
$ws.Range("O5").Value = -10
$ws.Range("P5").Value = 69
$ws.Range("Q5").Value = -40
$ws.Range("R5").Value = 6
$ws.Range("S5").Value = -1
$ws.Range("T5").Value = 9
$ws.Range("U5").Value = -49
$ws.Range("V5").Value = 99
$ws.Range("W5").Value = 0.23
$ws.Range("X5").Value = 0.76
$ws.Range("Y5").Value = 0.82
$ws.Range("Z5").Value = 0.97
$ws.Range("AA5").Value = 26.55
$ws.Range("AB5").Value = 1851.19
$ws.Range("AC5").Value = 157
$ws.Range("AD5").Value = 49.54
$ws.Range("AE5").Value = 20123
$ws.Range("AF5").Value = 0.39
$ws.Range("AG5").Value = 200
$ws.Range("AH5").Value = 2.57
$ws.Range("AI5").Value = 120.28
$ws.Range("AJ5").Value = 6867945

# Row 6 (2018/12) -- note: J6 and O6 stay empty (not present before either)
$ws.Range("D6").Value = 2061
$ws.Range("E6").Value = -39
$ws.Range("F6").Value = -39
$ws.Range("G6").Value = -56
$ws.Range("H6").Value = -67
$ws.Range("I6").Value = -63
$ws.Range("K6").Value = 1551
$ws.Range("L6").Value = 347
$ws.Range("M6").Value = 1204
$ws.Range("N6").Value = 1218
$ws.Range("P6").Value = 69
$ws.Range("Q6").Value = 9
$ws.Range("R6").Value = 15
$ws.Range("S6").Value = -13
$ws.Range("T6").Value = 10
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 97
$ws.Range("W6").Value = -1.9
$ws.Range("X6").Value = -3.25
$ws.Range("Y6").Value = -5
$ws.Range("Z6").Value = -4.2
$ws.Range("AA6").Value = 28.81
$ws.Range("AB6").Value = 1721.55
$ws.Range("AC6").Value = -917
$ws.Range("AD6").Value = -6.72
$ws.Range("AE6").Value = 18862
$ws.Range("AF6").Value = 0.33
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 3.25
$ws.Range("AI6").Value = -14.42
$ws.Range("AJ6").Value = 6867945

# Rows 7, 8, 9 (2019/12(E), 2020/12(E), 2021/12(E)):
# all numeric/forecast columns D:AI are cleared, leaving only A, B, C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
